# Apply updated Betfair Back/Lay odds for 2026-01-07 (rows 2-14)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.38
$ws.Range("G2").Value = 1.39
$ws.Range("H2").Value = 11
$ws.Range("J2").Value = 5.2
$ws.Range("K2").Value = 5.4
$ws.Range("P2").Value = 1.85
$ws.Range("Q2").Value = 2.12
$ws.Range("S2").Value = 3.95
$ws.Range("T2").Value = 2.6
$ws.Range("U2").Value = 1.59
$ws.Range("X2").Value = 14.5
$ws.Range("Y2").Value = 28
$ws.Range("Z2").Value = 130
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 48
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 44
$ws.Range("AI2").Value = 290
$ws.Range("AJ2").Value = 11
$ws.Range("AK2").Value = 18.5
$ws.Range("AL2").Value = 65
$ws.Range("AN2").Value = 8
# Row 3
$ws.Range("G3").Value = 2.9
$ws.Range("N3").Value = 3.8
$ws.Range("P3").Value = 1.92
$ws.Range("T3").Value = 1.78
$ws.Range("Z3").Value = 18.5
$ws.Range("AC3").Value = 7.4
# Row 4
$ws.Range("H4").Value = 8.199999999999999
$ws.Range("I4").Value = 8.4
$ws.Range("R4").Value = 1.81
$ws.Range("T4").Value = 1.7
$ws.Range("U4").Value = 2.36
$ws.Range("AB4").Value = 13
# Row 5
$ws.Range("F5").Value = 2.18
$ws.Range("G5").Value = 2.2
$ws.Range("P5").Value = 2.24
$ws.Range("Q5").Value = 1.78
$ws.Range("R5").Value = 1.47
# Row 6
$ws.Range("F6").Value = 1.83
$ws.Range("G6").Value = 1.84
$ws.Range("H6").Value = 5.3
$ws.Range("O6").Value = 1.39
$ws.Range("U6").Value = 1.9
$ws.Range("Z6").Value = 42
$ws.Range("AE6").Value = 85
# Row 7
$ws.Range("H7").Value = 2.42
$ws.Range("O7").Value = 1.32
$ws.Range("S7").Value = 3.55
$ws.Range("T7").Value = 1.77
$ws.Range("AG7").Value = 14.5
$ws.Range("AL7").Value = 48
# Row 8
$ws.Range("F8").Value = 1.91
$ws.Range("G8").Value = 1.92
$ws.Range("I8").Value = 4.8
$ws.Range("O8").Value = 1.37
$ws.Range("Q8").Value = 2.1
# Row 9
$ws.Range("F9").Value = 3.5
$ws.Range("H9").Value = 2.18
$ws.Range("I9").Value = 2.2
$ws.Range("J9").Value = 3.75
$ws.Range("Y9").Value = 11
$ws.Range("Z9").Value = 15
# Row 10
$ws.Range("F10").Value = 2.54
$ws.Range("G10").Value = 2.56
$ws.Range("H10").Value = 3.35
$ws.Range("I10").Value = 3.4
$ws.Range("J10").Value = 3.2
$ws.Range("K10").Value = 3.25
$ws.Range("O10").Value = 1.45
$ws.Range("Y10").Value = 10.5
$ws.Range("Z10").Value = 23
$ws.Range("AB10").Value = 8.800000000000001
$ws.Range("AD10").Value = 15
$ws.Range("AF10").Value = 15
$ws.Range("AG10").Value = 13
$ws.Range("AI10").Value = 65
$ws.Range("AJ10").Value = 36
$ws.Range("AN10").Value = 32
$ws.Range("AO10").Value = 55
# Row 11
$ws.Range("F11").Value = 2.38
$ws.Range("G11").Value = 2.4
$ws.Range("H11").Value = 3.65
$ws.Range("I11").Value = 3.7
$ws.Range("O11").Value = 1.47
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 1.93
$ws.Range("Z11").Value = 25
$ws.Range("AA11").Value = 75
$ws.Range("AC11").Value = 7
$ws.Range("AF11").Value = 14.5
$ws.Range("AH11").Value = 20
$ws.Range("AN11").Value = 29
$ws.Range("AO11").Value = 65
# Row 12
$ws.Range("F12").Value = 9.199999999999999
$ws.Range("G12").Value = 9.4
$ws.Range("H12").Value = 1.41
$ws.Range("I12").Value = 1.42
$ws.Range("J12").Value = 5.4
$ws.Range("K12").Value = 5.5
$ws.Range("O12").Value = 1.23
$ws.Range("P12").Value = 2.38
$ws.Range("Q12").Value = 1.69
$ws.Range("S12").Value = 2.74
$ws.Range("T12").Value = 1.97
$ws.Range("U12").Value = 1.99
$ws.Range("X12").Value = 22
$ws.Range("Y12").Value = 9.199999999999999
$ws.Range("Z12").Value = 8.800000000000001
$ws.Range("AF12").Value = 1000
$ws.Range("AO12").Value = 5.8
# Row 13
$ws.Range("F13").Value = 1.72
$ws.Range("G13").Value = 1.74
$ws.Range("H13").Value = 5.6
$ws.Range("K13").Value = 4.1
$ws.Range("Q13").Value = 1.88
$ws.Range("Y13").Value = 20
$ws.Range("Z13").Value = 48
$ws.Range("AF13").Value = 10.5
$ws.Range("AN13").Value = 10
# Row 14
$ws.Range("F14").Value = 5.1
$ws.Range("H14").Value = 1.76
$ws.Range("I14").Value = 1.77
$ws.Range("R14").Value = 1.46
$ws.Range("X14").Value = 18
$ws.Range("AA14").Value = 18.5
$ws.Range("AB14").Value = 20
$ws.Range("AE14").Value = 17.5
$ws.Range("AF14").Value = 42
$ws.Range("AH14").Value = 19
$ws.Range("AI14").Value = 32
$ws.Range("AJ14").Value = 150
$ws.Range("AL14").Value = 70
$ws.Range("AN14").Value = 70
